$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2003.7037
$ws.Range("J125").Value = 1845.6364
$ws.Range("L125").Value = 16610.7276
$ws.Range("N125").Value = -21530.7276
$ws.Range("H131").Value = 4525.326
$ws.Range("I131").Value = 416.66666
$ws.Range("J131").Value = 4811.9766
$ws.Range("K131").Value = 1249.99998
$ws.Range("L131").Value = 14435.9298
$ws.Range("M131").Value = 3790.00002
$ws.Range("N131").Value = -24515.9298

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 9335.571
$ws.Range("I38").Value = 2616.6667
$ws.Range("J38").Value = 14374.75
$ws.Range("K38").Value = 2616.6667
$ws.Range("L38").Value = 14374.75
$ws.Range("M38").Value = -2149.6667
$ws.Range("N38").Value = -15308.75
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H122").Value = 2344.875
$ws.Range("I122").Value = 2100.2
$ws.Range("J122").Value = 2752.6667
$ws.Range("K122").Value = 6300.599999999999
$ws.Range("L122").Value = 8258.000100000001
$ws.Range("M122").Value = -3850.599999999999
$ws.Range("N122").Value = -13158.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43523900
$ws.Range("J20").Value = 200001310
$ws.Range("L20").Value = 200001310
$ws.Range("N20").Value = -200001804
$ws.Range("H94").Value = 582.2917
$ws.Range("I94").Value = 579.73334
$ws.Range("J94").Value = 586.55554
$ws.Range("K94").Value = 579.73334
$ws.Range("L94").Value = 586.55554
$ws.Range("M94").Value = -128.73334
$ws.Range("N94").Value = -1488.55554
$ws.Range("H99").Value = 1722.4193
$ws.Range("I99").Value = 1127.7273
$ws.Range("K99").Value = 1127.7273
$ws.Range("M99").Value = 370.2727
$ws.Range("H134").Value = 2895.2258
$ws.Range("I134").Value = 2442.6667
$ws.Range("K134").Value = 7328.000100000001
$ws.Range("M134").Value = -4793.000100000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 5500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 5500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 5500
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -6254
$ws.Range("H46").Value = 5500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5922
$ws.Range("H68").Value = 16814.27
$ws.Range("J68").Value = 16814.27
$ws.Range("L68").Value = 16814.27
$ws.Range("N68").Value = -18312.27
$ws.Range("H71").Value = 16814.27
$ws.Range("J71").Value = 16814.27
$ws.Range("L71").Value = 50442.81
$ws.Range("N71").Value = -57930.81
$ws.Range("H122").Value = 988.8889
$ws.Range("I122").Value = 988.8889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2966.6667
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -516.6667000000002
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4725
$ws.Range("I132").Value = 6145.8887
$ws.Range("J132").Value = 2593.6667
$ws.Range("K132").Value = 18437.6661
$ws.Range("L132").Value = 7781.000100000001
$ws.Range("M132").Value = -15907.6661
$ws.Range("N132").Value = -12841.0001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 311.17648
$ws.Range("I11").Value = 307.27274
$ws.Range("J11").Value = 318.33334
$ws.Range("K11").Value = 921.81822
$ws.Range("L11").Value = 955.0000200000001
$ws.Range("M11").Value = -781.81822
$ws.Range("N11").Value = -1235.00002
$ws.Range("H34").Value = 803.26666
$ws.Range("I34").Value = 63
$ws.Range("J34").Value = 1072.4546
$ws.Range("K34").Value = 189
$ws.Range("L34").Value = 3217.3638
$ws.Range("M34").Value = -105
$ws.Range("N34").Value = -3385.3638
$ws.Range("H40").Value = 683
$ws.Range("I40").Value = 580.3333
$ws.Range("J40").Value = 760
$ws.Range("K40").Value = 2321.3332
$ws.Range("L40").Value = 3040
$ws.Range("M40").Value = -2252.3332
$ws.Range("N40").Value = -3178
$ws.Range("H46").Value = 127001
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 252002
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 756006
$ws.Range("M46").Value = -5909
$ws.Range("N46").Value = -756188
$ws.Range("H63").Value = 1700
$ws.Range("I63").Value = 1700
$ws.Range("K63").Value = 5100
$ws.Range("M63").Value = -4351
$ws.Range("H66").Value = 1700
$ws.Range("I66").Value = 1700
$ws.Range("K66").Value = 15300
$ws.Range("M66").Value = -11556
$ws.Range("H113").Value = 759.13513
$ws.Range("I113").Value = 1704.25
$ws.Range("J113").Value = 498.4138
$ws.Range("K113").Value = 5112.75
$ws.Range("L113").Value = 1495.2414
$ws.Range("M113").Value = -2942.75
$ws.Range("N113").Value = -5835.2414
$ws.Range("H131").Value = 773.23
$ws.Range("J131").Value = 785.81445
$ws.Range("L131").Value = 2357.44335
$ws.Range("N131").Value = -12437.44335
$ws.Range("H137").Value = 4597441
$ws.Range("I137").Value = 112622.89
$ws.Range("J137").Value = 7702315
$ws.Range("K137").Value = 337868.67
$ws.Range("L137").Value = 23106945
$ws.Range("M137").Value = -332768.67
$ws.Range("N137").Value = -23117145

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1641.5454
$ws.Range("I113").Value = 1309
$ws.Range("J113").Value = 1918.6666
$ws.Range("K113").Value = 1309
$ws.Range("L113").Value = 1918.6666
$ws.Range("M113").Value = 861
$ws.Range("N113").Value = -6258.6666
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4320.074
$ws.Range("I7").Value = 4294.2354
$ws.Range("J7").Value = 4364
$ws.Range("K7").Value = 4294.2354
$ws.Range("L7").Value = 4364
$ws.Range("M7").Value = -4182.2354
$ws.Range("N7").Value = -4588
$ws.Range("H16").Value = 424479.16
$ws.Range("I16").Value = 84331.914
$ws.Range("J16").Value = 716033.9399999999
$ws.Range("K16").Value = 84331.914
$ws.Range("L16").Value = 716033.9399999999
$ws.Range("M16").Value = -84161.914
$ws.Range("N16").Value = -716373.9399999999
$ws.Range("H32").Value = 1700
$ws.Range("J32").Value = 2400
$ws.Range("L32").Value = 2400
$ws.Range("N32").Value = -3034
$ws.Range("H40").Value = 85599.914
$ws.Range("I40").Value = 334969.66
$ws.Range("J40").Value = 2476.6667
$ws.Range("K40").Value = 334969.66
$ws.Range("L40").Value = 2476.6667
$ws.Range("M40").Value = -334833.66
$ws.Range("N40").Value = -2748.6667
$ws.Range("H61").Value = 1767.5883
$ws.Range("I61").Value = 1736.4
$ws.Range("J61").Value = 1812.1428
$ws.Range("K61").Value = 1736.4
$ws.Range("L61").Value = 1812.1428
$ws.Range("M61").Value = -1534.4
$ws.Range("N61").Value = -2216.1428
$ws.Range("H113").Value = 1767.5883
$ws.Range("I113").Value = 1736.4
$ws.Range("J113").Value = 1812.1428
$ws.Range("K113").Value = 1736.4
$ws.Range("L113").Value = 1812.1428
$ws.Range("M113").Value = 433.5999999999999
$ws.Range("N113").Value = -6152.1428
$ws.Range("H123").Value = 24751.8
$ws.Range("J123").Value = 31659.666
$ws.Range("L123").Value = 31659.666
$ws.Range("N123").Value = -41459.666
$ws.Range("H126").Value = 4320.074
$ws.Range("I126").Value = 4294.2354
$ws.Range("J126").Value = 4364
$ws.Range("K126").Value = 12882.7062
$ws.Range("L126").Value = 13092
$ws.Range("M126").Value = -10412.7062
$ws.Range("N126").Value = -18032
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 201059.9
$ws.Range("I107").Value = 1228.5714
$ws.Range("J107").Value = 667333
$ws.Range("K107").Value = 3685.7142
$ws.Range("L107").Value = 2001999
$ws.Range("M107").Value = -1765.7142
$ws.Range("N107").Value = -2005839
$ws.Range("H122").Value = 2531.5
$ws.Range("I122").Value = 1966.6666
$ws.Range("J122").Value = 2773.5715
$ws.Range("K122").Value = 5899.9998
$ws.Range("L122").Value = 8320.7145
$ws.Range("M122").Value = -3449.9998
$ws.Range("N122").Value = -13220.7145
$ws.Range("H135").Value = 32328.666
$ws.Range("J135").Value = 32328.666
$ws.Range("L135").Value = 32328.666
$ws.Range("N135").Value = -42468.666
$ws.Range("H137").Value = 30000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 30000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 30000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -40200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
